# Swap the values of columns D and E (codeforiati:group-code / codeforiati:group-name,
# and their corresponding per-row code/name values) for every used row in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
